$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("aragon"): dimension identifier changes to the curated sdmx dimension,
# and its representation changes from skos:Concept to a URI-based code list.
$ws.Range("B2").Value = "sdmx-dimension:refArea"
$ws.Range("B4").Value = "URI-Comunidad"

# Column C ("situacion-profesional"): now curated as a measure instead of a dimension.
$ws.Range("C2").Value = "iaest-measure:situacion-profesional"
$ws.Range("C3").Value = "medida"
$ws.Range("C4").Value = "xsd:int"

# Column D ("sexo"): now curated as a measure instead of a dimension.
$ws.Range("D2").Value = "iaest-measure:sexo"
$ws.Range("D3").Value = "medida"
$ws.Range("D4").Value = "xsd:int"

# Row 5 held the per-dimension mapping workbook references; no longer needed.
$ws.Rows.Item(5).Delete()
